# Add the LiveCampus timetable paste instructions to Sheet1!A1, matching
# the "Add check of student_intensive_lectures' validation" commit:
#  - A1 gets the Japanese instruction text (wrapped)
#  - Column A is narrowed, row 1 is made tall enough for the wrapped text
#  - The sheet is set up for A4 portrait printing

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ここに、LiveCampusからコピーした個人時間割を貼り付けてください。"
$ws.Range("A1").WrapText = $true

# ColumnWidth is specified in characters of the workbook's default font; the
# engine stores the width in "characters + 5px padding" units, so asking for
# 15 literally round-trips to 15.8333 on disk. Back the padding out so the
# persisted <col width="..."/> lands on exactly 15.
$ws.Columns(1).ColumnWidth = 85 / 6

# Row height is stored in points 1:1.
$ws.Rows(1).RowHeight = 90

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
